# Plantilla Lista de Tareas de la Iteración - apply edits described in commit:
#  - Mark CU 10/12 mockup (row16), CU 09/11 (row... wait kept for reference) as done
#  - Update descripcion CU 06/08 (row17) and CU 05/07 (row18) status + hours
#  - Update row 15 estimate / consumed hours
#  - Add 4 new task rows (19-22) for mockups CU 14/16, CU 13/15 and
#    descripciones CU 09/11, CU 10/12
#  - Update frozen-pane / selection view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Existing rows 15-18: update estimated/consumed hours and status
# ---------------------------------------------------------------------------

# Row 15 (Realizar mockup de CU 09 y 11): estimate 2h, consumed 3h on day 7
$ws.Range("G15").Value = 2
$ws.Range("AC15").Value = 3

# Row 16 (Realizar mockup de CU 10 y 12): done, estimate 2h
$ws.Range("F16").Value = "Hecho"
$ws.Range("G16").Value = 2
$ws.Range("AC16").Value = 2

# Row 17 (Realizar descripción de CU 06 y 08): done, estimate 1h
$ws.Range("F17").Value = "Hecho"
$ws.Range("G17").Value = 1
$ws.Range("AC17").Value = 1

# Row 18 (Realizar descripción de CU 05 y 07): done, estimate 1h
$ws.Range("F18").Value = "Hecho"
$ws.Range("G18").Value = 1
$ws.Range("AC18").Value = 1

# ---------------------------------------------------------------------------
# 2. Add new rows 19-22 with the same layout/formatting as row 18
# ---------------------------------------------------------------------------

$ws.Range("B18:BA18").Copy()
$ws.Range("B19:BA22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 48
$ws.Rows.Item(22).RowHeight = 45.75

# Row 19: Realizar mockup de CU 14 y 16 (Mario, por iniciar, sin estimado)
$ws.Range("D19").Value = "Realizar mockup de CU 14 y 16"
$ws.Range("E19").Value = "Mario"
$ws.Range("F19").Value = "Por iniciar"

# Row 20: Realizar mockup de CU 13 y 15 (Victor, por iniciar, estimado 2h)
$ws.Range("D20").Value = "Realizar mockup de CU 13 y 15"
$ws.Range("E20").Value = "Victor"
$ws.Range("F20").Value = "Por iniciar"
$ws.Range("G20").Value = 2

# Row 21: Realizar descripciones de CU  09 y 11 (Mario, por iniciar, sin estimado)
$ws.Range("D21").Value = "Realizar descripciones de CU  09 y 11"
$ws.Range("E21").Value = "Mario"
$ws.Range("F21").Value = "Por iniciar"

# Row 22: Realizar descripciones de CU 10 y 12 (Victor, por iniciar, estimado 1h)
$ws.Range("D22").Value = "Realizar descripciones de CU 10 y 12"
$ws.Range("E22").Value = "Victor"
$ws.Range("F22").Value = "Por iniciar"
$ws.Range("G22").Value = 1

# Restante/consumo helper formulas mirrored from row 18's pattern, for each
# of the new rows (I, L, O, R, U, X, AA, AD, AG, AJ, AM, AP, AS, AV, AY are the
# "restante" columns; AZ/BA reproduce row 18's broken #REF! totals formula).
$restCols = @(
  @("I", "G", "H"),
  @("L", "I", "K"),
  @("O", "L", "N"),
  @("R", "O", "Q"),
  @("U", "R", "T"),
  @("X", "U", "W"),
  @("AA", "X", "Z"),
  @("AD", "AA", "AC"),
  @("AG", "AD", "AF"),
  @("AJ", "AG", "AI"),
  @("AM", "AJ", "AL"),
  @("AP", "AM", "AO"),
  @("AS", "AP", "AR"),
  @("AV", "AS", "AU"),
  @("AY", "AV", "AX")
)

foreach ($r in 19..22) {
  foreach ($triplet in $restCols) {
    $target = $triplet[0]
    $left = $triplet[1]
    $right = $triplet[2]
    $ws.Range("$target$r").Formula = "=$left$r-$right$r"
  }
  $azFormula = "=H$r+K$r+N$r+Q$r+T$r+W$r+Z$r+AC$r+AF$r+AI$r+AL$r+AO$r+AR$r+AU$r+AX$r+#REF!+#REF!+#REF!+#REF!+#REF!"
  $ws.Range("AZ$r").Formula = $azFormula
  $ws.Range("BA$r").Formula = "=G$r-AZ$r"
}

# ---------------------------------------------------------------------------
# 3. View state: scroll frozen panes down and move the active selection
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 7
$ws.Range("G24").Select()

Write-Host "Edit applied."
